$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Borders for the merged header cells ---
# Sheet1 C1 -> top+bottom border (built fresh from default style)
$c1 = $ws1.Range("C1")
$c1.Style = "Normal"
$c1.Borders.Item(8).LineStyle = 1
$c1.Borders.Item(9).LineStyle = 1

# Sheet1 D1 -> top+bottom+right border (built fresh: top, then right, then bottom,
# so every intermediate combination matches a border already present in the file)
$d1 = $ws1.Range("D1")
$d1.Style = "Normal"
$d1.Borders.Item(8).LineStyle = 1
$d1.Borders.Item(10).LineStyle = 1
$d1.Borders.Item(9).LineStyle = 1

# Propagate the same formatting to sheet2's equivalent header cells
$c1.Copy()
$ws2.Range("C1").PasteSpecial(-4122)
$ws2.Range("F1").PasteSpecial(-4122)

$d1.Copy()
$ws2.Range("D1").PasteSpecial(-4122)
$ws2.Range("G1").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Anonymize "fedcore" -> "approach" ---
$ws1.Range("C2").Value = "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# --- Clear the stray empty inline-string cell ---
$ws2.Range("G5").ClearContents()
